# Update countries & provincias Spain
#
# 1. Bump the "last updated" timestamp in A1.
# 2. Haiti's case counts grew enough that it now outranks Venezuela and
#    Guinea-Bisau in the (descending, by total cases) sort, so its row
#    moves from position 106 up to position 104 (right after Maldivas),
#    pushing Venezuela and Guinea-Bisau down one row each. Their own data
#    is untouched - only their row position shifts - while Haiti gets
#    fresh totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Timestamp refresh
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 04:05"

# 2) Re-sorted rows 104-106 (country column + the 7 data columns)
# Row 104: now Haiti, with updated counts
$ws.Range("A104").Value = "Haiti"
$ws.Range("B104").Value = 1320
$ws.Range("C104").Value = 146
$ws.Range("D104").Value = 22
$ws.Range("E104").Value = 1264
$ws.Range("F104").Value = 0
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 34

# Row 105: now Venezuela (unchanged data, shifted down from old row 104)
$ws.Range("A105").Value = "Venezuela"
$ws.Range("B105").Value = 1245
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 302
$ws.Range("E105").Value = 932
$ws.Range("F105").Value = 0
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 11

# Row 106: now Guinea-Bisau (unchanged data, shifted down from old row 105)
$ws.Range("A106").Value = "Guinea-Bisau"
$ws.Range("B106").Value = 1195
$ws.Range("C106").Value = 0
$ws.Range("D106").Value = 42
$ws.Range("E106").Value = 1146
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 7
